$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Free up the "ImageFile" header text so that, once every new stimulus
# filename below has been inserted into the shared-string table, re-adding
# "ImageFile" appends it after them instead of reusing its old slot.
$ws.Range("A1").Value = ""

# New neutral-stimuli rows (the old pictures were removed from the trial
# list and these new ones -- plus the accompanying "readme" reorganisation
# implied by the commit message -- were added instead).
$ws.Range("A5").Value = "Stimuli/140.jpg"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "Stimuli/143.jpg"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "Stimuli/7000.jpg"
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = "Stimuli/7002.jpg"
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = "Stimuli/7004.jpg"
$ws.Range("B9").Value = 1

$ws.Range("A10").Value = "Stimuli/7006.jpg"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "Stimuli/7009.jpg"
$ws.Range("B11").Value = 1

$ws.Range("A12").Value = "Stimuli/7021.jpg"
$ws.Range("B12").Value = 1

$ws.Range("A13").Value = "Stimuli/7025.jpg"
$ws.Range("B13").Value = 1

$ws.Range("A14").Value = "Stimuli/7041.jpg"
$ws.Range("B14").Value = 1

$ws.Range("A15").Value = "Stimuli/7100.jpg"
$ws.Range("B15").Value = 1

$ws.Range("A16").Value = "Stimuli/7150.jpg"
$ws.Range("B16").Value = 1

$ws.Range("A17").Value = "Stimuli/7185.jpg"
$ws.Range("B17").Value = 1

$ws.Range("A18").Value = "Stimuli/7211.jpg"
$ws.Range("B18").Value = 1

$ws.Range("A19").Value = "Stimuli/7224.jpg"
$ws.Range("B19").Value = 1

$ws.Range("A20").Value = "Stimuli/7233.jpg"
$ws.Range("B20").Value = 1

$ws.Range("A21").Value = "Stimuli/7235.jpg"
$ws.Range("B21").Value = 1

# Restore the header text now that all the new stimuli strings exist.
$ws.Range("A1").Value = "ImageFile"

# Match the saved selection state (whole-row selection on row 18).
$ws.Rows(18).Select() | Out-Null
